$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("contacts")

# Update the "firstname" column values so each row has a unique value
$ws.Range("A3").Value = "Testing2"
$ws.Range("A4").Value = "Testing3"
$ws.Range("A2").Value = "Testing4"

# Move the selection from column C to cell A2
$ws.Range("A2").Select() | Out-Null
